$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 97.59999999999999
$ws.Range("I39").Value = 97.59999999999999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 292.8
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 3.200000000000045
$ws.Range("N39").ClearContents()
$ws.Range("H43").Value = 4953.857
$ws.Range("J43").Value = 5446.1665
$ws.Range("L43").Value = 5446.1665
$ws.Range("N43").Value = -5584.1665
$ws.Range("H76").Value = 7987.375
$ws.Range("J76").Value = 7987.375
$ws.Range("L76").Value = 7987.375
$ws.Range("N76").Value = -8617.375
$ws.Range("H79").Value = 7987.375
$ws.Range("J79").Value = 7987.375
$ws.Range("L79").Value = 7987.375
$ws.Range("N79").Value = -10171.375
$ws.Range("H87").Value = 80000.5
$ws.Range("J87").Value = 80000.5
$ws.Range("L87").Value = 80000.5
$ws.Range("N87").Value = -82496.5
$ws.Range("H90").Value = 80000.5
$ws.Range("J90").Value = 80000.5
$ws.Range("L90").Value = 240001.5
$ws.Range("N90").Value = -252481.5
$ws.Range("H100").Value = 565.55
$ws.Range("I100").Value = 488.8125
$ws.Range("K100").Value = 488.8125
$ws.Range("M100").Value = 52.1875
$ws.Range("H127").Value = 557
$ws.Range("I127").Value = 557
$ws.Range("K127").Value = 1671
$ws.Range("M127").Value = 3289
$ws.Range("H137").Value = 3316.9092
$ws.Range("I137").Value = 2148.375
$ws.Range("J137").Value = 3984.6428
$ws.Range("K137").Value = 6445.125
$ws.Range("L137").Value = 11953.9284
$ws.Range("M137").Value = -3895.125
$ws.Range("N137").Value = -17053.9284
$ws.Range("H138").Value = 13233.214
$ws.Range("J138").Value = 13308.223
$ws.Range("L138").Value = 39924.669
$ws.Range("N138").Value = -50204.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14742.117
$ws.Range("I32").Value = 6406.543
$ws.Range("J32").Value = 23582.879
$ws.Range("K32").Value = 6406.543
$ws.Range("L32").Value = 23582.879
$ws.Range("M32").Value = -6119.543
$ws.Range("N32").Value = -24156.879
$ws.Range("H61").Value = 1488.5333
$ws.Range("I61").Value = 1440.6154
$ws.Range("K61").Value = 1440.6154
$ws.Range("M61").Value = -1228.6154
$ws.Range("H132").Value = 2033.762
$ws.Range("I132").Value = 1985.8049
$ws.Range("K132").Value = 5957.4147
$ws.Range("M132").Value = -3427.4147
$ws.Range("H133").Value = 124494
$ws.Range("J133").Value = 124494
$ws.Range("L133").Value = 124494
$ws.Range("N133").Value = -129554
$ws.Range("H136").Value = 1488.5333
$ws.Range("I136").Value = 1440.6154
$ws.Range("K136").Value = 4321.8462
$ws.Range("M136").Value = -1771.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 674
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -946
$ws.Range("H82").Value = 36249.75
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 99999
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 99999
$ws.Range("M82").Value = -14617
$ws.Range("N82").Value = -100765
$ws.Range("H85").Value = 36249.75
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 99999
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 99999
$ws.Range("M85").Value = -13674
$ws.Range("N85").Value = -102651
$ws.Range("H86").Value = 1455.7778
$ws.Range("I86").Value = 1486.1428
$ws.Range("J86").Value = 1349.5
$ws.Range("K86").Value = 1486.1428
$ws.Range("L86").Value = 1349.5
$ws.Range("M86").Value = -363.1428000000001
$ws.Range("N86").Value = -3595.5
$ws.Range("H89").Value = 1455.7778
$ws.Range("I89").Value = 1486.1428
$ws.Range("J89").Value = 1349.5
$ws.Range("K89").Value = 7430.714
$ws.Range("L89").Value = 6747.5
$ws.Range("M89").Value = -1814.714
$ws.Range("N89").Value = -17979.5
$ws.Range("H94").Value = 564.5
$ws.Range("I94").Value = 461.76923
$ws.Range("K94").Value = 461.76923
$ws.Range("M94").Value = -10.76922999999999
$ws.Range("H105").Value = 3257.3125
$ws.Range("I105").Value = 3736.2222
$ws.Range("K105").Value = 3736.2222
$ws.Range("M105").Value = -1989.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4615.4443
$ws.Range("I31").Value = 2206.6667
$ws.Range("K31").Value = 2206.6667
$ws.Range("M31").Value = -1911.6667
$ws.Range("H34").Value = 4615.4443
$ws.Range("I34").Value = 2206.6667
$ws.Range("K34").Value = 2206.6667
$ws.Range("M34").Value = -2004.6667
$ws.Range("H58").Value = 6598.75
$ws.Range("J58").Value = 7461
$ws.Range("L58").Value = 7461
$ws.Range("N58").Value = -7867
$ws.Range("H136").Value = 6598.75
$ws.Range("J136").Value = 7461
$ws.Range("L136").Value = 22383
$ws.Range("N136").Value = -27483
$ws.Range("H141").Value = 115000
$ws.Range("J141").Value = 115000
$ws.Range("L141").Value = 115000
$ws.Range("N141").Value = -125360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 472
$ws.Range("I110").Value = 472
$ws.Range("K110").Value = 1416
$ws.Range("M110").Value = 2674
$ws.Range("H114").Value = 650.2727
$ws.Range("I114").Value = 496.7143
$ws.Range("J114").Value = 919
$ws.Range("K114").Value = 1490.1429
$ws.Range("L114").Value = 2757
$ws.Range("M114").Value = 1763.8571
$ws.Range("N114").Value = -9265

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2999.6667
$ws.Range("I43").Value = 999.5
$ws.Range("K43").Value = 999.5
$ws.Range("M43").Value = -848.5
$ws.Range("H97").Value = 2449.9524
$ws.Range("I97").Value = 2665.4
$ws.Range("K97").Value = 2665.4
$ws.Range("M97").Value = -2169.4
$ws.Range("H132").Value = 5924.75
$ws.Range("I132").Value = 4937
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 14811
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -12281
$ws.Range("N132").Value = -31724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5271.727
$ws.Range("I61").Value = 5355.2856
$ws.Range("J61").Value = 5125.5
$ws.Range("K61").Value = 5355.2856
$ws.Range("L61").Value = 5125.5
$ws.Range("M61").Value = -5153.2856
$ws.Range("N61").Value = -5529.5
$ws.Range("H113").Value = 5271.727
$ws.Range("I113").Value = 5355.2856
$ws.Range("J113").Value = 5125.5
$ws.Range("K113").Value = 5355.2856
$ws.Range("L113").Value = 5125.5
$ws.Range("M113").Value = -3185.2856
$ws.Range("N113").Value = -9465.5
$ws.Range("H132").Value = 4482
$ws.Range("I132").Value = 2941.6155
$ws.Range("K132").Value = 8824.8465
$ws.Range("M132").Value = -6294.8465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1456.0526
$ws.Range("I122").Value = 1461.5625
$ws.Range("J122").Value = 1426.6666
$ws.Range("K122").Value = 4384.6875
$ws.Range("L122").Value = 4279.9998
$ws.Range("M122").Value = -1934.6875
$ws.Range("N122").Value = -9179.9998
$ws.Range("H132").Value = 1795.25
$ws.Range("I132").Value = 1282.6666
$ws.Range("J132").Value = 3333
$ws.Range("K132").Value = 3847.9998
$ws.Range("L132").Value = 9999
$ws.Range("M132").Value = -1317.9998
$ws.Range("N132").Value = -15059
